# Convert the M2Doc "m:if" complex field (fldChar begin / instrText* / fldChar
# end) into plain literal text runs "{m:if self.name = 'anydsl'}", mirroring
# the TokenIteratorFieldRewriterSplit parser change.

$d = $word.ActiveDocument

# Locate the field (there is exactly one in this template) and remember
# where its paragraph starts before removing it. (Field.Result.Paragraphs
# is unreliable for an empty field result, so find the enclosing paragraph
# by scanning the document's paragraphs instead.)
$field = $d.Fields.Item(1)
$codeStart = $field.Code.Start

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if (($codeStart -ge $p.Range.Start) -and ($codeStart -le $p.Range.End)) {
        $targetParagraph = $p
    }
}
$insertPos = $targetParagraph.Range.Start

# Delete the whole field (fldChar begin, every instrText run, fldChar end).
$field.Delete()

# The literal text, split the same way the rewritten template splits it into
# separate <w:t> runs.
$pieces = @("{m:if ", "self.name ", "=", " ", "'", "anydsl", "'}")

# Insert each piece as its own run. A transient bookmark is dropped between
# consecutive insertions so the engine doesn't coalesce neighbouring text
# into a single run; the bookmarks are removed again once all the text is
# in place, leaving only the plain <w:r><w:t>…</w:t></w:r> runs behind.
$pos = $insertPos
$tempBookmarks = @()
$i = 0
foreach ($piece in $pieces) {
    if ($i -gt 0) {
        $bookmarkName = "m2docSplit" + $i
        $d.Bookmarks.Add($bookmarkName, $d.Range($pos, $pos))
        $tempBookmarks += $bookmarkName
    }
    $d.Range($pos, $pos).InsertAfter($piece)
    $pos = $pos + $piece.Length
    $i = $i + 1
}

foreach ($bookmarkName in $tempBookmarks) {
    $d.Bookmarks($bookmarkName).Delete()
}
